$wb = $excel.ActiveWorkbook

# --- "demand" sheet: insert a new "location" column (C) before "item" ---
$wsDemand = $wb.Worksheets.Item("demand")

# Insert a new column at C, shifting item/due/status/... one column to the right
$wsDemand.Columns.Item(3).Insert()

# Populate the new "location" column: header + "NMI" for every data row
$wsDemand.Cells.Item(1, 3).Value = "location"
$wsDemand.Cells.Item(2, 3).Value = "NMI"
$wsDemand.Cells.Item(3, 3).Value = "NMI"
$wsDemand.Cells.Item(4, 3).Value = "NMI"
$wsDemand.Cells.Item(5, 3).Value = "NMI"

# --- "location" sheet: update the lingering selection left over from editing ---
$wsLocation = $wb.Worksheets.Item("location")
$wsLocation.Range("A2").Select()

# --- "demand" sheet becomes the active / selected tab, with B2 selected ---
$wsDemand.Range("B2").Select()
